# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets.
# Both sheets hold identical duplicated data, so the same row/value
# updates are applied to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 3117
    4  = 2885
    8  = 1544
    10 = 75
    12 = 1279
    14 = 411
    16 = 61
    21 = 2886
    22 = 347
    23 = 12
    24 = 65
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
